# Apply the "gh-pages" deploy update to the FHIR StructureDefinition workbook.
#
# Sheet "Metadata" (sheet1):
#   - Version value 5.0.0 -> 6.0.0
#   - Date value updated
#   - Publisher value set to "Alvearie Team"
#   - The duplicated "Contact" / "No display for ContactDetail" row is removed,
#     and the remaining row becomes "Jurisdiction" / "United States of America"
#
# Sheet "Elements" (sheet2):
#   - Root "Extension" element row: Short -> "Offset End",
#     Definition -> the StructureDefinition description text

$wb = $excel.ActiveWorkbook

$metadata = $wb.Worksheets.Item("Metadata")
$elements = $wb.Worksheets.Item("Elements")

# --- Metadata sheet -------------------------------------------------

# Version
$metadata.Range("B3").Value = "6.0.0"

# Date
$metadata.Range("B8").Value = "2022-01-21T20:46:54+00:00"

# Publisher now has a value
$metadata.Range("B9").Value = "Alvearie Team"

# Remove the duplicate "Contact" row (row 11); row 10 keeps "Contact" /
# "No display for ContactDetail" for now and gets overwritten below.
$metadata.Rows.Item(11).Delete()

# Row 10 becomes the "Jurisdiction" property
$metadata.Range("A10").Value = "Jurisdiction"
$metadata.Range("B10").Value = "United States of America"

# --- Elements sheet --------------------------------------------------

# Root Extension row (row 2): Short / Definition columns (K / L)
$elements.Range("K2").Value = "Offset End"
$elements.Range("L2").Value = "Offset location of the last character for the span of covered text in relation to the overall reference where this span of text appears"

Write-Host "edit complete"
